# Update loading_percent values for rows 2-25 (columns B, C, D, E, G, I, J)
# as per the "case with 380 kV done" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 24.40962925897914
$ws.Range("C2").Value = 27.2694669272689
$ws.Range("D2").Value = 15.4332243392616
$ws.Range("E2").Value = 16.89904869702617
$ws.Range("G2").Value = 3.796599285071753
$ws.Range("I2").Value = 51.27908416197352
$ws.Range("J2").Value = 9.692096973219188

# Row 3
$ws.Range("B3").Value = 24.28329111718783
$ws.Range("C3").Value = 26.77821748385764
$ws.Range("D3").Value = 15.40108689930717
$ws.Range("E3").Value = 16.86960474898793
$ws.Range("G3").Value = 3.803576977178254
$ws.Range("I3").Value = 49.98904320153159
$ws.Range("J3").Value = 9.714116602872037

# Row 4
$ws.Range("B4").Value = 24.21909728476391
$ws.Range("C4").Value = 26.48256505077505
$ws.Range("D4").Value = 15.3855095132459
$ws.Range("E4").Value = 16.85608198951321
$ws.Range("G4").Value = 3.808061544129464
$ws.Range("I4").Value = 49.18262102414104
$ws.Range("J4").Value = 9.729410298342394

# Row 5
$ws.Range("B5").Value = 24.19630947948678
$ws.Range("C5").Value = 26.36374224586579
$ws.Range("D5").Value = 15.3802047034992
$ws.Range("E5").Value = 16.85171421484761
$ws.Range("G5").Value = 3.809939712252839
$ws.Range("I5").Value = 48.85070979802597
$ws.Range("J5").Value = 9.736087269543162

# Row 6
$ws.Range("B6").Value = 24.19272930646391
$ws.Range("C6").Value = 26.34411641834777
$ws.Range("D6").Value = 15.37938677548318
$ws.Range("E6").Value = 16.85105787557402
$ws.Range("G6").Value = 3.810254650346639
$ws.Range("I6").Value = 48.79540718883226
$ws.Range("J6").Value = 9.73722279384363

# Row 7
$ws.Range("B7").Value = 24.2187763037075
$ws.Range("C7").Value = 26.48095564910398
$ws.Range("D7").Value = 15.38543375021104
$ws.Range("E7").Value = 16.8560184611121
$ws.Range("G7").Value = 3.808086668174747
$ws.Range("I7").Value = 49.17815763495062
$ws.Range("J7").Value = 9.729498547595748

# Row 8
$ws.Range("B8").Value = 24.36329795703747
$ws.Range("C8").Value = 27.09893745509143
$ws.Range("D8").Value = 15.42128017639434
$ws.Range("E8").Value = 16.88794944682085
$ws.Range("G8").Value = 3.798963828069077
$ws.Range("I8").Value = 50.83743702427954
$ws.Range("J8").Value = 9.699320454728841

# Row 9
$ws.Range("B9").Value = 24.75217608209525
$ws.Range("C9").Value = 28.35132170433358
$ws.Range("D9").Value = 15.52458356431682
$ws.Range("E9").Value = 16.98680036364543
$ws.Range("G9").Value = 3.782648058977753
$ws.Range("I9").Value = 53.96458563790581
$ws.Range("J9").Value = 9.654272930182071

# Row 10
$ws.Range("B10").Value = 25.10079865734526
$ws.Range("C10").Value = 29.28694068226288
$ws.Range("D10").Value = 15.6206039913485
$ws.Range("E10").Value = 17.08157698416907
$ws.Range("G10").Value = 3.771599744275859
$ws.Range("I10").Value = 56.16886532937227
$ws.Range("J10").Value = 9.629873922522671

# Row 11
$ws.Range("B11").Value = 25.27262755957872
$ws.Range("C11").Value = 29.71403829811704
$ws.Range("D11").Value = 15.66863690400328
$ws.Range("E11").Value = 17.1294967714891
$ws.Range("G11").Value = 3.766772922987149
$ws.Range("I11").Value = 57.14841843121663
$ws.Range("J11").Value = 9.620680408961709

# Row 12
$ws.Range("B12").Value = 25.33955368487326
$ws.Range("C12").Value = 29.87582918101919
$ws.Range("D12").Value = 15.6874484032871
$ws.Range("E12").Value = 17.14833149268383
$ws.Range("G12").Value = 3.764973401362433
$ws.Range("I12").Value = 57.51579294831423
$ws.Range("J12").Value = 9.617474527794457

# Row 13
$ws.Range("B13").Value = 25.32505807718677
$ws.Range("C13").Value = 29.84098441963361
$ws.Range("D13").Value = 15.68336938955608
$ws.Range("E13").Value = 17.14424452476101
$ws.Range("G13").Value = 3.765359707781465
$ws.Range("I13").Value = 57.4368340480348
$ws.Range("J13").Value = 9.61815269536994

# Row 14
$ws.Range("B14").Value = 25.27809664840844
$ws.Range("C14").Value = 29.72734862570429
$ws.Range("D14").Value = 15.67017208479846
$ws.Range("E14").Value = 17.13103253344934
$ws.Range("G14").Value = 3.766624310241648
$ws.Range("I14").Value = 57.17871514958315
$ws.Range("J14").Value = 9.620411129343973

# Row 15
$ws.Range("B15").Value = 25.24957198073377
$ws.Range("C15").Value = 29.65774653321264
$ws.Range("D15").Value = 15.66216930247879
$ws.Range("E15").Value = 17.12302938605567
$ws.Range("G15").Value = 3.767402589859365
$ws.Range("I15").Value = 57.02013973429801
$ws.Range("J15").Value = 9.621830407289336

# Row 16
$ws.Range("B16").Value = 25.08983186236338
$ws.Range("C16").Value = 29.25904527769729
$ws.Range("D16").Value = 15.6175522417926
$ws.Range("E16").Value = 17.07854182021098
$ws.Range("G16").Value = 3.771919165975288
$ws.Range("I16").Value = 56.10436282125379
$ws.Range("J16").Value = 9.630513213366925

# Row 17
$ws.Range("B17").Value = 24.995196535387
$ws.Range("C17").Value = 29.01472213904715
$ws.Range("D17").Value = 15.59129379438879
$ws.Range("E17").Value = 17.05247945682859
$ws.Range("G17").Value = 3.774740695139124
$ws.Range("I17").Value = 55.5364571683281
$ws.Range("J17").Value = 9.636328976523084

# Row 18
$ws.Range("B18").Value = 24.94201396237554
$ws.Range("C18").Value = 28.87434211352613
$ws.Range("D18").Value = 15.57660059384865
$ws.Range("E18").Value = 17.03794159261503
$ws.Range("G18").Value = 3.776382326171413
$ws.Range("I18").Value = 55.20764287099493
$ws.Range("J18").Value = 9.639853337813831

# Row 19
$ws.Range("B19").Value = 24.9242231023077
$ws.Range("C19").Value = 28.82684210760819
$ws.Range("D19").Value = 15.57169623065225
$ws.Range("E19").Value = 17.03309710228081
$ws.Range("G19").Value = 3.776941386956254
$ws.Range("I19").Value = 55.09594655377422
$ws.Range("J19").Value = 9.641077374193253

# Row 20
$ws.Range("B20").Value = 25.00514165695355
$ws.Range("C20").Value = 29.04071655705039
$ws.Range("D20").Value = 15.59404664336813
$ws.Range("E20").Value = 17.05520702104004
$ws.Range("G20").Value = 3.7744383991979
$ws.Range("I20").Value = 55.59713798901746
$ws.Range("J20").Value = 9.635691312353298

# Row 21
$ws.Range("B21").Value = 25.29184031318051
$ws.Range("C21").Value = 29.76072583105355
$ws.Range("D21").Value = 15.67403159101307
$ws.Range("E21").Value = 17.13489455313729
$ws.Range("G21").Value = 3.766252100780274
$ws.Range("I21").Value = 57.25462927987491
$ws.Range("J21").Value = 9.619740283839009

# Row 22
$ws.Range("B22").Value = 25.49002127169376
$ws.Range("C22").Value = 30.23156626155448
$ws.Range("D22").Value = 15.72993245011198
$ws.Range("E22").Value = 17.19098649904795
$ws.Range("G22").Value = 3.761066620638223
$ws.Range("I22").Value = 58.31704843749268
$ws.Range("J22").Value = 9.610921881180392

# Row 23
$ws.Range("B23").Value = 25.3832758236353
$ws.Range("C23").Value = 29.98029426583269
$ws.Range("D23").Value = 15.69976668439731
$ws.Range("E23").Value = 17.1606831691788
$ws.Range("G23").Value = 3.763819249259755
$ws.Range("I23").Value = 57.7519926433387
$ws.Range("J23").Value = 9.615480950351047

# Row 24
$ws.Range("B24").Value = 25.00064164994434
$ws.Range("C24").Value = 29.02896420536134
$ws.Range("D24").Value = 15.59280082443176
$ws.Range("E24").Value = 17.0539725007155
$ws.Range("G24").Value = 3.774575006440131
$ws.Range("I24").Value = 55.56971139187888
$ws.Range("J24").Value = 9.635979037278794

# Row 25
$ws.Range("B25").Value = 24.63579956198494
$ws.Range("C25").Value = 28.00913820306786
$ws.Range("D25").Value = 15.49309966878194
$ws.Range("E25").Value = 16.95616925036248
$ws.Range("G25").Value = 3.786895519182751
$ws.Range("I25").Value = 53.13390058311573
$ws.Range("J25").Value = 9.664938051289941

